# Apply the "2016-12-18 -> 2016-12-19" date refresh plus the two text
# fixes described in the commit ("Fixing hyphenation on 'best practices'
# per @jurph" + removing a stray duplicate endParaRPr in the notes).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" placeholder text from
#    2016-12-18 to 2016-12-19 everywhere it is cached: the slide master,
#    the handout master, the notes master, and every slide layout.
# ---------------------------------------------------------------------

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "2016-12-18") {
                    $tr.Text = "2016-12-19"
                }
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster
Update-DatePlaceholder $p.HandoutMaster
Update-DatePlaceholder $p.NotesMaster

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    Update-DatePlaceholder $p.SlideMaster.CustomLayouts.Item($j)
}

# ---------------------------------------------------------------------
# 2) Slide 2, "TextBox 18": fix hyphenation, "best-practices" ->
#    "best practices".
# ---------------------------------------------------------------------

$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $sh = $slide2.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tf = $sh.TextFrame
        if ($tf.HasText) {
            $tr = $tf.TextRange
            if ($tr.Text -like "*best-practices*") {
                $tr.Text = [string]::Concat([char]0x201C, "I detect hygiene issues and operator activity that does not follow best practices.", [char]0x201D)
            }
        }
    }
}

